# Add a new "canonical SMILES" column (column D) to the microstate list sheet.
# Column D mirrors column C ("canonical isomeric SMILES") but with the explicit
# cis/trans bond-stereo markers ("/" and "\") stripped out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("D2").Value = "canonical SMILES"

# Data rows 3-14 hold one microstate each; column C already has the
# canonical isomeric SMILES string, so derive column D from it.
for ($r = 3; $r -le 14; $r++) {
    $isomericSmiles = $ws.Range("C$r").Text
    $canonicalSmiles = $isomericSmiles.Replace("/", "").Replace("\", "")
    $ws.Range("D$r").Value = $canonicalSmiles
}

# Give the new column a sensible width, matching the other SMILES columns.
$ws.Columns.Item(4).ColumnWidth = 36
